$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update existing data cells (B/C columns) with the new text values.
#    Re-using the existing hyperlink cells (C2:C4) keeps their r:id intact.
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "Poznań I"
$ws.Range("B3").Value = "Poznań I"
$ws.Range("B4").Value = "Bieszczady I"
$ws.Range("C2").Value = "Wp.pl"
$ws.Range("C3").Value = "onet.pl"
$ws.Range("C4").Value = "okopress"

# ---------------------------------------------------------------------------
# 2. Extend the table down to row 33 with blank (formatted) rows.
# ---------------------------------------------------------------------------
$ws.Range("A5:C33").Value = ""

# ---------------------------------------------------------------------------
# 3. Formatting: thin borders around the whole data block, then a date
#    number format for column A.
# ---------------------------------------------------------------------------
$ws.Range("A2:C33").Borders.LineStyle = 1
$ws.Range("A2:A33").NumberFormat = "yyyy\-mm\-dd;@"

# ---------------------------------------------------------------------------
# 4. Column widths.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 21
$ws.Columns.Item(2).ColumnWidth = 29.6
$ws.Columns.Item(3).ColumnWidth = 66.3

# ---------------------------------------------------------------------------
# 5. Data validation (dropdown list) on B2:B33.
# ---------------------------------------------------------------------------
$ws.Range("B2:B33").Validation.Add(3, 1, 1, '"Van,PC,Poznań I,Poznań II, Bieszczady I, Bieszczady II"') | Out-Null

# ---------------------------------------------------------------------------
# 6. Selection / active cell, matching the saved view state.
# ---------------------------------------------------------------------------
$ws.Range("K11").Select() | Out-Null
